$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 74: PEL.NS
$ws.Cells.Item(74, 1).Value = "PEL.NS"
$ws.Cells.Item(74, 2).Value = 37803
$ws.Cells.Item(74, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(74, 3).Value = 30
$ws.Cells.Item(74, 4).Value = 37073
$ws.Cells.Item(74, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(74, 5).Value = 22
$ws.Cells.Item(74, 6).Value = 21.93608474731445
$ws.Cells.Item(74, 7).Value = 37438
$ws.Cells.Item(74, 7).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(74, 8).Value = 26
$ws.Cells.Item(74, 9).Value = 21.79438591003418
$ws.Cells.Item(74, 10).Value = "Low"
$ws.Cells.Item(74, 11).Value = -0.03542470932006836
$ws.Cells.Item(74, 12).Value = 22.71542835235596
$ws.Cells.Item(74, 13).Value = 3
$ws.Cells.Item(74, 14).Value = 1
$ws.Cells.Item(74, 15).Value = 2

# Row 75: PEL.NS
$ws.Cells.Item(75, 1).Value = "PEL.NS"
$ws.Cells.Item(75, 2).Value = 44835
$ws.Cells.Item(75, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(75, 3).Value = 107
$ws.Cells.Item(75, 4).Value = 42826
$ws.Cells.Item(75, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(75, 5).Value = 85
$ws.Cells.Item(75, 6).Value = 1741.661987304688
$ws.Cells.Item(75, 7).Value = 44470
$ws.Cells.Item(75, 7).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(75, 8).Value = 103
$ws.Cells.Item(75, 9).Value = 1751.69140625
$ws.Cells.Item(75, 10).Value = "High"
$ws.Cells.Item(75, 11).Value = 0.55718994140625
$ws.Cells.Item(75, 12).Value = 1694.300842285156
$ws.Cells.Item(75, 13).Value = 3
$ws.Cells.Item(75, 14).Value = 1
$ws.Cells.Item(75, 15).Value = 2
